$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 'SEC_001, SEC_002, SEC_004, SEC_005, SEC_006, SEC_011, SEC_012, SEC_014, SEC_015, SEC_016, SEC_018, SEC_019, SEC_021, SEC_023'
$ws.Range("B5").Value = 'SEC_000, SEC_003, SEC_007, SEC_008, SEC_009, SEC_010, SEC_013, SEC_017, SEC_020, SEC_022, SEC_024, SEC_025, SEC_026, SEC_027'
$ws.Range("B6").Value = 'SEC_001, SEC_002, SEC_003, SEC_005, SEC_007, SEC_008, SEC_010, SEC_011, SEC_012, SEC_013, SEC_019, SEC_024, SEC_025, SEC_027'
$ws.Range("B7").Value = 'SEC_000, SEC_004, SEC_006, SEC_009, SEC_014, SEC_015, SEC_016, SEC_017, SEC_018, SEC_020, SEC_021, SEC_022, SEC_023, SEC_026'
$ws.Range("B8").Value = 'SEC_000, SEC_001, SEC_002, SEC_005, SEC_009, SEC_012, SEC_014, SEC_016, SEC_018, SEC_019, SEC_020, SEC_022, SEC_025, SEC_026'
$ws.Range("B9").Value = 'SEC_003, SEC_004, SEC_006, SEC_007, SEC_008, SEC_010, SEC_011, SEC_013, SEC_015, SEC_017, SEC_021, SEC_023, SEC_024, SEC_027'
$ws.Range("B10").Value = 'SEC_001, SEC_002, SEC_003, SEC_006, SEC_009, SEC_013, SEC_016, SEC_017, SEC_020, SEC_021, SEC_024, SEC_025, SEC_026, SEC_027'
$ws.Range("B11").Value = 'SEC_000, SEC_004, SEC_005, SEC_007, SEC_008, SEC_010, SEC_011, SEC_012, SEC_014, SEC_015, SEC_018, SEC_019, SEC_022, SEC_023'
$ws.Range("B12").Value = 'SEC_000, SEC_004, SEC_006, SEC_007, SEC_008, SEC_009, SEC_010, SEC_011, SEC_013, SEC_017, SEC_020, SEC_021, SEC_022, SEC_027'
$ws.Range("B13").Value = 'SEC_001, SEC_002, SEC_003, SEC_005, SEC_012, SEC_014, SEC_015, SEC_016, SEC_018, SEC_019, SEC_023, SEC_024, SEC_025, SEC_026'
$ws.Range("G14").Value = '(4, 2), (4, 3)'
$ws.Range("G15").Value = '(2, 0), (2, 1)'
$ws.Range("G16").Value = '(4, 0), (4, 1)'
$ws.Range("G17").Value = '(0, 0), (0, 1)'
$ws.Range("G18").Value = '(2, 4), (2, 5)'
$ws.Range("G19").Value = '(1, 0), (1, 1)'
$ws.Range("G20").Value = '(3, 4), (3, 5)'
$ws.Range("G21").Value = '(2, 0), (2, 1)'
$ws.Range("G22").Value = '(0, 2), (0, 3)'
$ws.Range("G23").Value = '(3, 4), (3, 5)'
$ws.Range("G24").Value = '(1, 6), (1, 7)'
$ws.Range("G25").Value = '(2, 6), (2, 7)'
$ws.Range("G26").Value = '(3, 2), (3, 3)'
$ws.Range("G27").Value = '(2, 2), (2, 3)'
$ws.Range("G28").Value = '(3, 0), (3, 1)'
$ws.Range("G29").Value = '(0, 4), (0, 5)'
$ws.Range("G30").Value = '(1, 6), (1, 7)'
$ws.Range("G31").Value = '(4, 0), (4, 1)'
$ws.Range("G32").Value = '(3, 6), (3, 7)'
$ws.Range("G33").Value = '(1, 6), (1, 7)'
$ws.Range("G34").Value = '(0, 0), (0, 1)'
$ws.Range("G35").Value = '(4, 0), (4, 1)'
$ws.Range("G36").Value = '(3, 4), (3, 5)'
$ws.Range("G37").Value = '(4, 4), (4, 5)'
$ws.Range("G38").Value = '(3, 4), (3, 5)'
$ws.Range("G39").Value = '(0, 4), (0, 5)'
$ws.Range("G40").Value = '(1, 2), (1, 3)'
$ws.Range("G41").Value = '(3, 6), (3, 7)'
$ws.Range("G42").Value = '(4, 6), (4, 7)'
$ws.Range("G43").Value = '(0, 0), (0, 1)'
$ws.Range("G44").Value = '(2, 0), (2, 1)'
$ws.Range("G45").Value = '(3, 6), (3, 7)'
$ws.Range("G46").Value = '(0, 6), (0, 7)'
$ws.Range("G47").Value = '(4, 4), (4, 5)'
$ws.Range("G48").Value = '(0, 0), (0, 1)'
$ws.Range("G49").Value = '(2, 2), (2, 3)'
$ws.Range("G50").Value = '(1, 6), (1, 7)'
$ws.Range("G51").Value = '(3, 2), (3, 3)'
$ws.Range("G52").Value = '(2, 6), (2, 7)'
$ws.Range("G53").Value = '(3, 6), (3, 7)'
$ws.Range("G54").Value = '(1, 0), (1, 1)'
$ws.Range("G55").Value = '(2, 0), (2, 1)'
$ws.Range("G56").Value = '(0, 0), (0, 1)'
$ws.Range("G57").Value = '(2, 4), (2, 5)'
$ws.Range("G58").Value = '(1, 6), (1, 7)'
$ws.Range("G59").Value = '(3, 6), (3, 7)'
$ws.Range("G60").Value = '(1, 0), (1, 1)'
$ws.Range("G62").Value = '(1, 0), (1, 1)'
$ws.Range("G63").Value = '(2, 0), (2, 1)'
$ws.Range("G64").Value = '(3, 6), (3, 7)'
$ws.Range("G65").Value = '(3, 0), (3, 1)'
$ws.Range("G66").Value = '(4, 4), (4, 5)'
$ws.Range("G67").Value = '(4, 6), (4, 7)'
$ws.Range("G68").Value = '(0, 4), (0, 5)'
$ws.Range("G69").Value = '(2, 4), (2, 5)'
$ws.Range("G70").Value = '(1, 0), (1, 1)'
$ws.Range("G72").Value = '(3, 4), (3, 5)'
$ws.Range("G73").Value = '(3, 2), (3, 3)'
$ws.Range("G74").Value = '(0, 0), (0, 1)'
$ws.Range("G76").Value = '(0, 2), (0, 3)'
$ws.Range("G79").Value = '(4, 4), (4, 5)'
$ws.Range("G80").Value = '(2, 0), (2, 1)'
$ws.Range("G81").Value = '(3, 4), (3, 5)'
$ws.Range("G82").Value = '(2, 0), (2, 1)'
$ws.Range("G83").Value = '(1, 6), (1, 7)'
$ws.Range("G84").Value = '(3, 2), (3, 3)'
$ws.Range("G85").Value = '(1, 4), (1, 5)'
$ws.Range("G86").Value = '(4, 4), (4, 5)'
$ws.Range("G87").Value = '(1, 0), (1, 1)'
$ws.Range("G88").Value = '(3, 4), (3, 5)'
$ws.Range("G89").Value = '(2, 0), (2, 1)'
$ws.Range("G90").Value = '(3, 0), (3, 1)'
$ws.Range("G91").Value = '(1, 0), (1, 1)'
$ws.Range("G93").Value = '(0, 2), (0, 3)'
$ws.Range("G94").Value = '(2, 6), (2, 7)'
$ws.Range("G95").Value = '(4, 6), (4, 7)'
$ws.Range("G96").Value = '(1, 6), (1, 7)'
$ws.Range("G97").Value = '(3, 6), (3, 7)'
$ws.Range("G98").Value = '(3, 4), (3, 5)'
$ws.Range("G99").Value = '(4, 6), (4, 7)'
